$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Đơn phụ phẫu 1" (2nd sheet)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

# Row 2: service group changed
$ws.Range("G2").Value = "Tiêm"

# Row 3: service group + payment figures changed
$ws.Range("G3").Value = "Tiểu phẫu"
$ws.Range("O3").Value = 500000
$ws.Range("P3").Value = 3500000
$ws.Range("Q3").Value = 500000
$ws.Range("Z3").Value = 350000

# Row 4: service group changed
$ws.Range("G4").Value = "Tiêm"

# Row 5: used to be the "Tổng" row - becomes a regular data row
$ws.Range("A5").Value = "HD-LUXURY"
$ws.Range("B5").Value = 542
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "07-13-2024"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").Value = "LONG XUYÊN"
$ws.Range("E5").Value = "Huỳnh thị bé sáu"
$ws.Range("F5").Value = "Cá nhân"
$ws.Range("G5").Value = "Tiểu phẫu"
$ws.Range("H5").Value = "Cắt mí"
$ws.Range("I5").Value = "Nguyễn Phúc Nam"
$ws.Range("J5").Value = 11500000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 11500000
$ws.Range("N5").Value = 9000000
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 9000000
$ws.Range("Q5").Value = 2500000
$ws.Range("R5").Value = "Lâm Thị Mỹ Hằng"
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = "Đào Vương Anh"
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 50000
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0.15
$ws.Range("Y5").Value = 0
$ws.Range("Z5").Value = 1350000
$ws.Range("AA5").Value = 0

# Row 6: new data row
$ws.Range("A6").Value = "HD-LUXURY"
$ws.Range("B6").Value = 551
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "07-14-2024"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").Value = "LONG XUYÊN"
$ws.Range("E6").Value = "Kim anh"
$ws.Range("F6").Value = "Cá nhân"
$ws.Range("G6").Value = "Nâng mũi"
$ws.Range("H6").Value = "Thu cánh mũi"
$ws.Range("I6").Value = "Nguyễn Phúc Nam"
$ws.Range("J6").Value = 4000000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 4000000
$ws.Range("N6").Value = 4000000
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 4000000
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = "Phạm Thanh Hoàng"
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = "Đào Vương Anh"
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 100000
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0.1
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 400000
$ws.Range("AA6").Value = 0

# Row 7: new "Tổng" row (totals for rows 2-6)
$ws.Range("A7").Value = "Tổng"
$ws.Range("B7").Value = 5
$ws.Range("J7").Value = 23850000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 23850000
$ws.Range("N7").Value = 20350000
$ws.Range("O7").Value = 500000
$ws.Range("P7").Value = 20850000
$ws.Range("Q7").Value = 3000000
$ws.Range("V7").Value = 200000
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0.55
$ws.Range("Y7").Value = 0
$ws.Range("Z7").Value = 2535000
$ws.Range("AA7").Value = 0

# ---------------------------------------------------------------------------
# Sheet "Lương" (3rd sheet)
# ---------------------------------------------------------------------------
$wsL = $wb.Worksheets.Item(3)

$wsL.Range("B2").Value = 16.5
$wsL.Range("B3").Value = 577500

# Insert a new row ("Ứng lương tại CẦN THƠ") after the CẦN THƠ block (old row 10),
# pushing everything below down by one.
$wsL.Rows.Item(11).Insert()
$wsL.Range("A11").Value = "Ứng lương tại CẦN THƠ"
$wsL.Range("B11").Value = 0

# Lương cơ bản tại LONG XUYÊN now gets a concrete value
$wsL.Range("B12").Value = 1767857.142857143

# Công phụ phẫu 1 tại LONG XUYÊN value changed
$wsL.Range("B17").Value = 200000

# Insert a new row ("Ứng lương tại LONG XUYÊN") after the LONG XUYÊN block
$wsL.Rows.Item(19).Insert()
$wsL.Range("A19").Value = "Ứng lương tại LONG XUYÊN"
$wsL.Range("B19").Value = 0

# Insert a new row ("Ứng lương tại SÓC TRĂNG") after the SÓC TRĂNG block (after old "Công phụ phẫu 2 tại SÓC TRĂNG")
$wsL.Rows.Item(27).Insert()
$wsL.Range("A27").Value = "Ứng lương tại SÓC TRĂNG"
$wsL.Range("B27").Value = 0

# New summary rows at the bottom
$wsL.Range("A28").Value = "Tổng lương tại CẦN THƠ"
$wsL.Range("B28").Value = 0
$wsL.Range("A29").Value = "Tổng lương tại LONG XUYÊN"
$wsL.Range("B29").Value = 6876607.142857143
$wsL.Range("A30").Value = "Tổng lương tại SÓC TRĂNG"
$wsL.Range("B30").Value = 0
$wsL.Range("A31").Value = "Tổng lương"
$wsL.Range("B31").Value = 6876607.142857143
